# Loan RBI, Variable Instalments
# Switch to the "Repayment schedule" sheet, insert a new (blank) column
# before column N, give it the same width as the column to its left (M),
# and leave the selection where the user ended up (L15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this the active sheet/tab (mirrors tabSelected moving from
# "Transactions" to "Repayment schedule", and workbook activeTab 3 -> 2).
$ws.Activate()

# Insert a blank column at N - shifts old N/O/P -> O/P/Q.
$ws.Columns("N").Insert()

# New column inherits the width of the column immediately to its left.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Final selection after the edit.
$ws.Range("L15").Select()
